$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet (was "US_Presidents Excel Tutorial Da")
$ws.Name = "US_Presidents"

# Update the active selection / active cell shown in the saved view
$ws.Range("C6").Select()
